$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level cell updates (price + 1h volume change refresh)
$cellUpdates = @{
    "D2" = "27.201.26"
    "E2" = "  +1.27%  "
    "D3" = "1.642.50"
    "E3" = "  +0.35%  "
    "E4" = "  +0.09%  "
    "D5" = "217.17"
    "E5" = "  +0.19%  "
    "D6" = "0.524"
    "E6" = "  +2.66%  "
    "E7" = "  +0.09%  "
    "E8" = "  -0.45%  "
    "E9" = "  +0.72%  "
    "D10" = "19.98"
    "E10" = "  +0.33%  "
    "D11" = "0.0849"
    "E11" = "  +0.47%  "
    "D12" = "1.871.72"
    "E12" = "  +0.32%  "
    "D13" = "1.637.50"
    "E13" = "  -0.06%  "
    "E14" = "  +0.76%  "
    "E15" = "  +2.63%  "
    "E16" = "  +0.69%  "
    "D17" = "27.191.94"
    "E17" = "  +1.23%  "
    "D18" = "0.0₃0741"
    "E18" = "  +1.75%  "
    "D19" = "218.91"
    "E19" = "  -0.17%  "
    "E20" = "  +0.16%  "
    "D21" = "6.97"
    "E21" = "  +3.43%  "
    "D22" = "4.42"
    "E22" = "  +0.80%  "
    "D23" = "2.51"
    "E23" = "  +3.24%  "
    "D24" = "9.12"
    "E24" = "  -0.32%  "
    "D25" = "147.42"
    "E25" = "  +0.26%  "
    "E26" = "  -0.14%  "
    "D27" = "7.45"
    "E27" = "  +1.31%  "
    "D28" = "0.120"
    "E28" = "  +0.47%  "
    "D29" = "15.73"
    "E29" = "  -0.29%  "
    "E30" = "  +1.07%  "
    "E31" = "  +0.45%  "
    "E32" = "  +1.56%  "
    "E33" = "  +0.34%  "
    "B34" = "LidoDAOToken"
    "C34" = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
    "D34" = "1.58"
    "E34" = "  +1.30%  "
    "B35" = "Maker"
    "C35" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
    "D35" = "1.303.56"
    "E35" = "  +3.73%  "
    "E36" = "  +1.52%  "
    "D37" = "0.0177"
    "E37" = "  -0.74%  "
    "E38" = "  +3.26%  "
    "D39" = "0.858"
    "E39" = "  +3.28%  "
    "E40" = "  +0.08%  "
    "E41" = "  +0.09%  "
    "E42" = "  +6.07%  "
    "E43" = "  -1.54%  "
    "D44" = "1.781.79"
    "E44" = "  +0.13%  "
    "D45" = "61.83"
    "E45" = "  +0.47%  "
    "D46" = "91.92"
    "E46" = "  +0.30%  "
    "E47" = "  +1.53%  "
    "E48" = "  +2.54%  "
    "E49" = "  -0.02%  "
    "D50" = "7.64"
    "E50" = "  +0.32%  "
    "D51" = "0.0965"
    "E51" = "  +0.64%  "
}

foreach ($addr in $cellUpdates.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cellUpdates[$addr]
    $rng.Style = "Normal"
}

